# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 79
$ws1.Range("F6").Value  = 41
$ws1.Range("F8").Value  = 122
$ws1.Range("F9").Value  = 8923
$ws1.Range("F13").Value = 1014
$ws1.Range("F18").Value = 295
$ws1.Range("F20").Value = 236
$ws1.Range("F21").Value = 1103

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 79
$ws4.Range("F7").Value  = 41
$ws4.Range("F10").Value = 122
$ws4.Range("F11").Value = 8923
$ws4.Range("F15").Value = 1014
$ws4.Range("F20").Value = 295
$ws4.Range("F22").Value = 236
$ws4.Range("F23").Value = 1103
